$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capitalize the three section header rows (gender / mother's education /
# wealth quintile) across Kyrgyz, Russian and English columns. Written in
# this specific order so the regenerated shared-string table lands in the
# same sequence as the published workbook.
$ws.Cells.Item(14, 1).Value = "Жынысы боюнча"
$ws.Cells.Item(14, 2).Value = "По полу"
$ws.Cells.Item(17, 2).Value = "Образование матери "
$ws.Cells.Item(17, 1).Value = "Энесинин билими "
$ws.Cells.Item(23, 2).Value = "Квинтиль по индексу благосостояния"
$ws.Cells.Item(14, 3).Value = "By sex"
$ws.Cells.Item(17, 3).Value = "Education of mother"
$ws.Cells.Item(23, 3).Value = "Wealth quintile"

# Restore the active selection to the sheet's default (A1) so the saved
# view state matches the republished workbook.
$ws.Range("A1").Select()
